# Update LR-pair TPM values (Tgfb1-Tgfbr1) with the newly computed numbers.
# Columns G,H,I,J,M,N,O,P,Q,R,S,T are replaced for rows 2-10; K and L stay unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column order for each row's data array below:
# G, H, I, J, M, N, O, P, Q, R, S, T
$newValues = @{
    2  = @(35.906979, 107.720937, 0.6107087147789413, 0.6107087147789412, 7.106976666666665, 21.32093, 0.1598176868560746, 0.1598176868560746, 255.19006192349, 2296.71055731141, 0.0976020541388166, 0.09760205413881658)
    3  = @(35.906979, 107.720937, 0.6107087147789413, 0.6107087147789412, 28.23943666666667, 84.71831, 0.6350325402576649, 0.6350325402576648, 1013.99285936183, 9125.935734256469, 0.3878199065035648, 0.3878199065035647)
    4  = @(35.906979, 107.720937, 0.6107087147789413, 0.6107087147789412, 9.122861666666667, 27.368585, 0.2051497728862606, 0.2051497728862606, 327.574402284905, 2948.169620564145, 0.1252867541365599, 0.1252867541365599)
    5  = @(17.04862266666667, 51.14586800000001, 0.2899643113254147, 0.2899643113254147, 7.106976666666665, 21.32093, 0.1598176868560746, 0.1598176868560746, 121.1641634908045, 1090.47747141724, 0.04634142550684246, 0.04634142550684243)
    6  = @(17.04862266666667, 51.14586800000001, 0.2899643113254147, 0.2899643113254147, 28.23943666666667, 84.71831, 0.6350325402576649, 0.6350325402576648, 481.4435000492312, 4332.99150044308, 0.1841367732050425, 0.1841367732050424)
    7  = @(17.04862266666667, 51.14586800000001, 0.2899643113254147, 0.2899643113254147, 9.122861666666667, 27.368585, 0.2051497728862606, 0.2051497728862606, 155.5322261951978, 1399.79003575678, 0.0594861126135298, 0.05948611261352978)
    8  = @(5.839988000000001, 17.519964, 0.09932697389564409, 0.09932697389564407, 7.106976666666665, 21.32093, 0.1598176868560746, 0.1598176868560746, 41.50465844961333, 373.54192604652, 0.01587420721041554, 0.01587420721041554)
    9  = @(5.839988000000001, 17.519964, 0.09932697389564409, 0.09932697389564407, 28.23943666666667, 84.71831, 0.6350325402576649, 0.6350325402576648, 164.9179712600934, 1484.26174134084, 0.06307586054905763, 0.06307586054905762)
    10 = @(5.839988000000001, 17.519964, 0.09932697389564409, 0.09932697389564407, 9.122861666666667, 27.368585, 0.2051497728862606, 0.2051497728862606, 53.27740265899335, 479.49662393094, 0.02037690613617092, 0.02037690613617092)
}

$cols = @("G", "H", "I", "J", "M", "N", "O", "P", "Q", "R", "S", "T")

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $vals[$i]
    }
}
